$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Cells.Item(2, 4) "63.190.61"
Set-TextValue $ws.Cells.Item(2, 5) "  +3.42%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.033.79"
Set-TextValue $ws.Cells.Item(3, 5) "  +1.84%  "
Set-TextValue $ws.Cells.Item(4, 5) "  -0.05%  "
Set-TextValue $ws.Cells.Item(5, 4) "595.61"
Set-TextValue $ws.Cells.Item(5, 5) "  +0.11%  "
Set-TextValue $ws.Cells.Item(6, 4) "154.28"
Set-TextValue $ws.Cells.Item(6, 5) "  +8.12%  "
Set-TextValue $ws.Cells.Item(7, 5) "  -0.02%  "
Set-TextValue $ws.Cells.Item(8, 4) "3.030.89"
Set-TextValue $ws.Cells.Item(8, 5) "  +1.84%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.516"
Set-TextValue $ws.Cells.Item(9, 5) "  +0.29%  "
Set-TextValue $ws.Cells.Item(10, 4) "6.96"
Set-TextValue $ws.Cells.Item(10, 5) "  +16.09%  "
Set-TextValue $ws.Cells.Item(11, 5) "  +3.43%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.465"
Set-TextValue $ws.Cells.Item(12, 5) "  +2.91%  "
Set-TextValue $ws.Cells.Item(13, 4) "0.0000235"
Set-TextValue $ws.Cells.Item(13, 5) "  +4.13%  "
Set-TextValue $ws.Cells.Item(14, 4) "35.71"
Set-TextValue $ws.Cells.Item(14, 5) "  +4.80%  "
Set-TextValue $ws.Cells.Item(15, 5) "  -0.52%  "
Set-TextValue $ws.Cells.Item(16, 4) "3.534.51"
Set-TextValue $ws.Cells.Item(16, 5) "  +1.81%  "
Set-TextValue $ws.Cells.Item(17, 4) "7.11"
Set-TextValue $ws.Cells.Item(17, 5) "  +3.88%  "
Set-TextValue $ws.Cells.Item(18, 4) "63.120.99"
Set-TextValue $ws.Cells.Item(18, 5) "  +3.08%  "
Set-TextValue $ws.Cells.Item(19, 4) "3.034.20"
Set-TextValue $ws.Cells.Item(19, 5) "  +1.94%  "
Set-TextValue $ws.Cells.Item(20, 4) "453.26"
Set-TextValue $ws.Cells.Item(20, 5) "  +1.14%  "
Set-TextValue $ws.Cells.Item(21, 4) "14.31"
Set-TextValue $ws.Cells.Item(21, 5) "  +2.24%  "
Set-TextValue $ws.Cells.Item(22, 4) "0.698"
Set-TextValue $ws.Cells.Item(22, 5) "  +2.57%  "
Set-TextValue $ws.Cells.Item(23, 4) "7.54"
Set-TextValue $ws.Cells.Item(23, 5) "  +3.51%  "
Set-TextValue $ws.Cells.Item(24, 4) "83.36"
Set-TextValue $ws.Cells.Item(24, 5) "  +1.88%  "
Set-TextValue $ws.Cells.Item(25, 4) "11.46"
Set-TextValue $ws.Cells.Item(25, 5) "  +9.86%  "
Set-TextValue $ws.Cells.Item(26, 4) "2.35"
Set-TextValue $ws.Cells.Item(26, 5) "  +9.00%  "
Set-TextValue $ws.Cells.Item(27, 4) "12.42"
Set-TextValue $ws.Cells.Item(27, 5) "  +4.52%  "
Set-TextValue $ws.Cells.Item(28, 5) "  -0.09%  "
Set-TextValue $ws.Cells.Item(29, 4) "7.62"
Set-TextValue $ws.Cells.Item(29, 5) "  +7.01%  "
Set-TextValue $ws.Cells.Item(30, 4) "2.29"
Set-TextValue $ws.Cells.Item(30, 5) "  +12.04%  "
Set-TextValue $ws.Cells.Item(31, 5) "  +1.67%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.999"
Set-TextValue $ws.Cells.Item(32, 5) "  -0.11%  "
Set-TextValue $ws.Cells.Item(33, 4) "27.66"
Set-TextValue $ws.Cells.Item(33, 5) "  +1.76%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.112"
Set-TextValue $ws.Cells.Item(34, 5) "  +3.43%  "
Set-TextValue $ws.Cells.Item(35, 4) "0.0₃0870"
Set-TextValue $ws.Cells.Item(35, 5) "  +8.02%  "
Set-TextValue $ws.Cells.Item(36, 5) "  +2.99%  "
Set-TextValue $ws.Cells.Item(37, 4) "5.92"
Set-TextValue $ws.Cells.Item(37, 5) "  +2.63%  "
Set-TextValue $ws.Cells.Item(38, 4) "3.22"
Set-TextValue $ws.Cells.Item(38, 5) "  +14.39%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.133"
Set-TextValue $ws.Cells.Item(39, 5) "  +10.30%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.12"
Set-TextValue $ws.Cells.Item(40, 5) "  +3.97%  "
Set-TextValue $ws.Cells.Item(41, 4) "50.56"
Set-TextValue $ws.Cells.Item(41, 5) "  +1.11%  "
Set-TextValue $ws.Cells.Item(42, 4) "9.12"
Set-TextValue $ws.Cells.Item(42, 5) "  +2.13%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.312"
Set-TextValue $ws.Cells.Item(43, 5) "  +16.56%  "
Set-TextValue $ws.Cells.Item(44, 4) "44.94"
Set-TextValue $ws.Cells.Item(44, 5) "  +17.62%  "
Set-TextValue $ws.Cells.Item(45, 4) "394.41"
Set-TextValue $ws.Cells.Item(45, 5) "  +2.66%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.0362"
Set-TextValue $ws.Cells.Item(46, 5) "  +3.99%  "
Set-TextValue $ws.Cells.Item(47, 4) "2.731.35"
Set-TextValue $ws.Cells.Item(47, 5) "  +1.32%  "
Set-TextValue $ws.Cells.Item(48, 4) "132.95"
Set-TextValue $ws.Cells.Item(48, 5) "  +2.71%  "
Set-TextValue $ws.Cells.Item(49, 4) "25.85"
Set-TextValue $ws.Cells.Item(49, 5) "  +11.87%  "
Set-TextValue $ws.Cells.Item(51, 5) "  +8.15%  "
